$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 385-404 (values shifted down by 3 rows; rows 385-387 get new entries) ---
$ws.Cells.Item(385,4).Value2 = 44615
$ws.Cells.Item(385,10).Value2 = 600
$ws.Cells.Item(385,11).Value2 = 9000
$ws.Cells.Item(385,12).Value2 = 9500
$ws.Cells.Item(385,13).Value2 = 9250
$ws.Cells.Item(385,14).Value2 = '$/bandeja 18 kilos'
$ws.Cells.Item(385,16).Value2 = 514
$ws.Cells.Item(385,17).Value2 = 18
$ws.Cells.Item(386,4).Value2 = 44615
$ws.Cells.Item(386,10).Value2 = 1000
$ws.Cells.Item(386,12).Value2 = 6500
$ws.Cells.Item(386,13).Value2 = 6250
$ws.Cells.Item(386,16).Value2 = 417
$ws.Cells.Item(387,4).Value2 = 44615
$ws.Cells.Item(387,9).Value2 = 'Segunda'
$ws.Cells.Item(387,10).Value2 = 600
$ws.Cells.Item(387,11).Value2 = 5000
$ws.Cells.Item(387,12).Value2 = 5500
$ws.Cells.Item(387,13).Value2 = 5250
$ws.Cells.Item(387,14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(387,16).Value2 = 350
$ws.Cells.Item(387,17).Value2 = 15
$ws.Cells.Item(388,4).Value2 = 44522
$ws.Cells.Item(388,9).Value2 = 'Primera'
$ws.Cells.Item(388,10).Value2 = 360
$ws.Cells.Item(388,11).Value2 = 7000
$ws.Cells.Item(388,12).Value2 = 8000
$ws.Cells.Item(388,13).Value2 = 7500
$ws.Cells.Item(388,14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(388,16).Value2 = 500
$ws.Cells.Item(388,17).Value2 = 15
$ws.Cells.Item(389,4).Value2 = 44543
$ws.Cells.Item(389,10).Value2 = 800
$ws.Cells.Item(389,11).Value2 = 6000
$ws.Cells.Item(389,12).Value2 = 7000
$ws.Cells.Item(389,13).Value2 = 6500
$ws.Cells.Item(389,16).Value2 = 433
$ws.Cells.Item(390,4).Value2 = 44167
$ws.Cells.Item(390,10).Value2 = 300
$ws.Cells.Item(390,11).Value2 = 4500
$ws.Cells.Item(390,12).Value2 = 5000
$ws.Cells.Item(390,13).Value2 = 4750
$ws.Cells.Item(390,14).Value2 = '$/caja 10 kilos'
$ws.Cells.Item(390,16).Value2 = 475
$ws.Cells.Item(390,17).Value2 = 10
$ws.Cells.Item(391,4).Value2 = 44167
$ws.Cells.Item(391,9).Value2 = 'Segunda'
$ws.Cells.Item(391,10).Value2 = 300
$ws.Cells.Item(391,11).Value2 = 3500
$ws.Cells.Item(391,12).Value2 = 4000
$ws.Cells.Item(391,13).Value2 = 3750
$ws.Cells.Item(391,15).Value2 = 'Región del Maule'
$ws.Cells.Item(391,16).Value2 = 375
$ws.Cells.Item(392,4).Value2 = 44277
$ws.Cells.Item(392,10).Value2 = 160
$ws.Cells.Item(392,11).Value2 = 4000
$ws.Cells.Item(392,12).Value2 = 4500
$ws.Cells.Item(392,13).Value2 = 4250
$ws.Cells.Item(392,14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(392,15).Value2 = 'Región del Maule'
$ws.Cells.Item(392,16).Value2 = 283
$ws.Cells.Item(392,17).Value2 = 15
$ws.Cells.Item(393,4).Value2 = 44258
$ws.Cells.Item(393,10).Value2 = 220
$ws.Cells.Item(393,11).Value2 = 8500
$ws.Cells.Item(393,12).Value2 = 9000
$ws.Cells.Item(393,13).Value2 = 8727
$ws.Cells.Item(393,14).Value2 = '$/bandeja 18 kilos'
$ws.Cells.Item(393,16).Value2 = 485
$ws.Cells.Item(393,17).Value2 = 18
$ws.Cells.Item(394,4).Value2 = 44390
$ws.Cells.Item(394,8).Value2 = 'Larga vida'
$ws.Cells.Item(394,11).Value2 = 6000
$ws.Cells.Item(394,12).Value2 = 6500
$ws.Cells.Item(394,13).Value2 = 6250
$ws.Cells.Item(394,14).Value2 = '$/caja 10 kilos'
$ws.Cells.Item(394,15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(394,16).Value2 = 625
$ws.Cells.Item(394,17).Value2 = 10
$ws.Cells.Item(395,4).Value2 = 44349
$ws.Cells.Item(395,11).Value2 = 5500
$ws.Cells.Item(395,12).Value2 = 6000
$ws.Cells.Item(395,13).Value2 = 5750
$ws.Cells.Item(395,16).Value2 = 575
$ws.Cells.Item(396,4).Value2 = 44285
$ws.Cells.Item(396,10).Value2 = 120
$ws.Cells.Item(396,11).Value2 = 5000
$ws.Cells.Item(396,12).Value2 = 5500
$ws.Cells.Item(396,13).Value2 = 5250
$ws.Cells.Item(396,16).Value2 = 350
$ws.Cells.Item(397,4).Value2 = 44285
$ws.Cells.Item(397,8).Value2 = 'Semiduro'
$ws.Cells.Item(397,11).Value2 = 4000
$ws.Cells.Item(397,12).Value2 = 4500
$ws.Cells.Item(397,13).Value2 = 4250
$ws.Cells.Item(397,14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(397,15).Value2 = 'Provincia de Diguillín'
$ws.Cells.Item(397,16).Value2 = 283
$ws.Cells.Item(397,17).Value2 = 15
$ws.Cells.Item(398,4).Value2 = 44498
$ws.Cells.Item(398,11).Value2 = 6000
$ws.Cells.Item(398,12).Value2 = 6500
$ws.Cells.Item(398,13).Value2 = 6250
$ws.Cells.Item(398,16).Value2 = 625
$ws.Cells.Item(399,4).Value2 = 44179
$ws.Cells.Item(399,10).Value2 = 300
$ws.Cells.Item(399,11).Value2 = 7500
$ws.Cells.Item(399,12).Value2 = 8000
$ws.Cells.Item(399,13).Value2 = 7750
$ws.Cells.Item(399,16).Value2 = 517
$ws.Cells.Item(400,4).Value2 = 44418
$ws.Cells.Item(400,9).Value2 = 'Primera'
$ws.Cells.Item(400,10).Value2 = 120
$ws.Cells.Item(400,11).Value2 = 5000
$ws.Cells.Item(400,12).Value2 = 5500
$ws.Cells.Item(400,13).Value2 = 5250
$ws.Cells.Item(400,14).Value2 = '$/caja 10 kilos'
$ws.Cells.Item(400,15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(400,16).Value2 = 525
$ws.Cells.Item(400,17).Value2 = 10
$ws.Cells.Item(401,4).Value2 = 44335
$ws.Cells.Item(401,10).Value2 = 300
$ws.Cells.Item(401,14).Value2 = '$/caja 10 kilos'
$ws.Cells.Item(401,15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(401,16).Value2 = 675
$ws.Cells.Item(401,17).Value2 = 10
$ws.Cells.Item(402,4).Value2 = 44552
$ws.Cells.Item(402,10).Value2 = 2000
$ws.Cells.Item(402,11).Value2 = 8000
$ws.Cells.Item(402,12).Value2 = 9000
$ws.Cells.Item(402,13).Value2 = 8500
$ws.Cells.Item(402,14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(402,15).Value2 = 'Región del Maule'
$ws.Cells.Item(402,16).Value2 = 567
$ws.Cells.Item(402,17).Value2 = 15
$ws.Cells.Item(403,4).Value2 = 44552
$ws.Cells.Item(403,9).Value2 = 'Segunda'
$ws.Cells.Item(403,10).Value2 = 600
$ws.Cells.Item(403,11).Value2 = 7000
$ws.Cells.Item(403,12).Value2 = 7000
$ws.Cells.Item(403,13).Value2 = 7000
$ws.Cells.Item(403,14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(403,15).Value2 = 'Región del Maule'
$ws.Cells.Item(403,16).Value2 = 467
$ws.Cells.Item(403,17).Value2 = 15
$ws.Cells.Item(404,4).Value2 = 44544
$ws.Cells.Item(404,10).Value2 = 2000
$ws.Cells.Item(404,11).Value2 = 6500
$ws.Cells.Item(404,12).Value2 = 7000
$ws.Cells.Item(404,13).Value2 = 6750
$ws.Cells.Item(404,15).Value2 = 'Región del Maule'
$ws.Cells.Item(404,16).Value2 = 450

# --- Append new rows 405-407 (values that were pushed off the front, i.e. old rows 385-387) ---
# Row 405
$ws.Cells.Item(405,1).Value2 = 7
$ws.Cells.Item(405,2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(405,3).Value2 = 'Ñuble'
$ws.Cells.Item(405,4).Value2 = 44160
$ws.Cells.Item(405,5).Value2 = 16
$ws.Cells.Item(405,6).Value2 = 100112020
$ws.Cells.Item(405,7).Value2 = 'Tomate'
$ws.Cells.Item(405,8).Value2 = 'Larga vida'
$ws.Cells.Item(405,9).Value2 = 'Primera'
$ws.Cells.Item(405,10).Value2 = 130
$ws.Cells.Item(405,11).Value2 = 11000
$ws.Cells.Item(405,12).Value2 = 12000
$ws.Cells.Item(405,13).Value2 = 11615
$ws.Cells.Item(405,14).Value2 = '$/bandeja 18 kilos'
$ws.Cells.Item(405,15).Value2 = 'Provincia de Talca'
$ws.Cells.Item(405,16).Value2 = 645
$ws.Cells.Item(405,17).Value2 = 18
$ws.Cells.Item(405,18).Value2 = 'Hortaliza'
$ws.Cells.Item(405,4).NumberFormat = $ws.Cells.Item(384,4).NumberFormat

# Row 406
$ws.Cells.Item(406,1).Value2 = 7
$ws.Cells.Item(406,2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(406,3).Value2 = 'Ñuble'
$ws.Cells.Item(406,4).Value2 = 44160
$ws.Cells.Item(406,5).Value2 = 16
$ws.Cells.Item(406,6).Value2 = 100112020
$ws.Cells.Item(406,7).Value2 = 'Tomate'
$ws.Cells.Item(406,8).Value2 = 'Larga vida'
$ws.Cells.Item(406,9).Value2 = 'Primera'
$ws.Cells.Item(406,10).Value2 = 280
$ws.Cells.Item(406,11).Value2 = 5000
$ws.Cells.Item(406,12).Value2 = 5500
$ws.Cells.Item(406,13).Value2 = 5268
$ws.Cells.Item(406,14).Value2 = '$/caja 10 kilos'
$ws.Cells.Item(406,15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(406,16).Value2 = 527
$ws.Cells.Item(406,17).Value2 = 10
$ws.Cells.Item(406,18).Value2 = 'Hortaliza'
$ws.Cells.Item(406,4).NumberFormat = $ws.Cells.Item(384,4).NumberFormat

# Row 407
$ws.Cells.Item(407,1).Value2 = 7
$ws.Cells.Item(407,2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(407,3).Value2 = 'Ñuble'
$ws.Cells.Item(407,4).Value2 = 44160
$ws.Cells.Item(407,5).Value2 = 16
$ws.Cells.Item(407,6).Value2 = 100112020
$ws.Cells.Item(407,7).Value2 = 'Tomate'
$ws.Cells.Item(407,8).Value2 = 'Larga vida'
$ws.Cells.Item(407,9).Value2 = 'Primera'
$ws.Cells.Item(407,10).Value2 = 290
$ws.Cells.Item(407,11).Value2 = 8000
$ws.Cells.Item(407,12).Value2 = 8500
$ws.Cells.Item(407,13).Value2 = 8241
$ws.Cells.Item(407,14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(407,15).Value2 = 'Provincia de Talca'
$ws.Cells.Item(407,16).Value2 = 549
$ws.Cells.Item(407,17).Value2 = 15
$ws.Cells.Item(407,18).Value2 = 'Hortaliza'
$ws.Cells.Item(407,4).NumberFormat = $ws.Cells.Item(384,4).NumberFormat

